# Re-arrange the template header row to add the new IRS-derived columns
# (Income, Deductibility Code, Asset Code, Assets, Income Code) while
# keeping Revenue and Mission Statement as the last two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (A1:K1):
#   A EIN | B Name | C Street | D City | E State | F ZIP | G GuideStar Link
#   H Link | I NTEE Code | J Revenue | K Mission Statement
#
# Target layout (A1:P1):
#   A EIN | B Name | C Street | D City | E State | F ZIP | G GuideStar Link
#   H Link | I NTEE Code | J Deductibility Code | K Asset Code | L Assets
#   M Income Code | N Income | O Revenue | P Mission Statement

# Move the existing Revenue / Mission Statement headers out to their
# final destination columns first, so they aren't clobbered.
$ws.Range("O1").Value = $ws.Range("J1").Text
$ws.Range("P1").Value = $ws.Range("K1").Text

# Write the new header cells in between. Order matters for how the
# shared-string table de-dupes new entries, so "Income" (N1) is written
# first, followed by the rest in column order.
$ws.Range("N1").Value = "Income"
$ws.Range("J1").Value = "Deductibility Code"
$ws.Range("K1").Value = "Asset Code"
$ws.Range("L1").Value = "Assets"
$ws.Range("M1").Value = "Income Code"

# Update column widths for the touched/new columns.
$ws.Columns.Item(10).ColumnWidth = 15.7265625
$ws.Columns.Item(11).ColumnWidth = 10
$ws.Columns.Item(13).ColumnWidth = 11.453125

# Update the view: scroll back to the top-left and move the selection.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D7").Select()
